$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores (matrix of p-values, symmetric with 1s on diagonal) ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.9773267524167553
$wsP.Range("D2").Value = 0.7149657756623857
$wsP.Range("E2").Value = 0.7208224960399416
$wsP.Range("F2").Value = 0.9216599800514955

$wsP.Range("B3").Value = 0.9773267524167553
$wsP.Range("D3").Value = 0.7387495294995581
$wsP.Range("E3").Value = 0.7031693339468537
$wsP.Range("F3").Value = 0.8970997383002497

$wsP.Range("B4").Value = 0.7149657756623857
$wsP.Range("C4").Value = 0.7387495294995581
$wsP.Range("E4").Value = 0.6365610379192925
$wsP.Range("F4").Value = 0.7954826854642367

$wsP.Range("B5").Value = 0.7208224960399416
$wsP.Range("C5").Value = 0.7031693339468537
$wsP.Range("D5").Value = 0.6365610379192925
$wsP.Range("F5").Value = 0.6218696113266748

$wsP.Range("B6").Value = 0.9216599800514955
$wsP.Range("C6").Value = 0.8970997383002497
$wsP.Range("D6").Value = 0.7954826854642367
$wsP.Range("E6").Value = 0.6218696113266748

# --- Sheet: Estadisticos_DM (anti-symmetric matrix of DM statistics) ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 0.02893258824943891
$wsE.Range("D2").Value = -0.3726831629691838
$wsE.Range("E2").Value = 0.364649136154689
$wsE.Range("F2").Value = -0.1001311393821628

$wsE.Range("B3").Value = -0.02893258824943891
$wsE.Range("D3").Value = -0.3402093738907146
$wsE.Range("E3").Value = 0.3889435192824905
$wsE.Range("F3").Value = -0.1316946066247819

$wsE.Range("B4").Value = 0.3726831629691838
$wsE.Range("C4").Value = 0.3402093738907146
$wsE.Range("E4").Value = 0.4830006637928517
$wsE.Range("F4").Value = 0.2641902202458279

$wsE.Range("B5").Value = -0.364649136154689
$wsE.Range("C5").Value = -0.3889435192824905
$wsE.Range("D5").Value = -0.4830006637928517
$wsE.Range("F5").Value = -0.5043356080234241

$wsE.Range("B6").Value = 0.1001311393821628
$wsE.Range("C6").Value = 0.1316946066247819
$wsE.Range("D6").Value = -0.2641902202458279
$wsE.Range("E6").Value = 0.5043356080234241
